# "Natmi following Dr Hou advice": split the myeloid sending cluster into
# M1/M2 subtypes (adds an "M1" sending-cluster block) and refresh the NATMI
# ligand/receptor specificity metrics for every sending/target cluster pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mmp2"
$ws.Range("C2").Value = "Sdc2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.333723
$ws.Range("H2").Value = 37.001169
$ws.Range("I2").Value = 0.03887240735130859
$ws.Range("J2").Value = 0.0388724073513086
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.248835333333334
$ws.Range("N2").Value = 6.746506
$ws.Range("O2").Value = 0.03590294220158827
$ws.Range("P2").Value = 0.03590294220158827
$ws.Range("Q2").Value = 27.73651207394601
$ws.Range("R2").Value = 249.628608665514
$ws.Range("S2").Value = 0.001395633794370628
$ws.Range("T2").Value = 0.001395633794370628

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mmp2"
$ws.Range("C3").Value = "Sdc2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.333723
$ws.Range("H3").Value = 37.001169
$ws.Range("I3").Value = 0.03887240735130859
$ws.Range("J3").Value = 0.0388724073513086
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 44.29005966666667
$ws.Range("N3").Value = 132.870179
$ws.Range("O3").Value = 0.7070964373190639
$ws.Range("P3").Value = 0.7070964373190639
$ws.Range("Q3").Value = 546.2613275821391
$ws.Range("R3").Value = 4916.351948239252
$ws.Range("S3").Value = 0.02748654074812569
$ws.Range("T3").Value = 0.0274865407481257

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mmp2"
$ws.Range("C4").Value = "Sdc2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.333723
$ws.Range("H4").Value = 37.001169
$ws.Range("I4").Value = 0.03887240735130859
$ws.Range("J4").Value = 0.0388724073513086
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 16.09762433333333
$ws.Range("N4").Value = 48.292873
$ws.Range("O4").Value = 0.2570006204793478
$ws.Range("P4").Value = 0.2570006204793479
$ws.Range("Q4").Value = 198.543639485393
$ws.Range("R4").Value = 1786.892755368537
$ws.Range("S4").Value = 0.009990232808812269
$ws.Range("T4").Value = 0.009990232808812272

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mmp2"
$ws.Range("C5").Value = "Sdc2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 301.5732116666666
$ws.Range("H5").Value = 904.7196349999999
$ws.Range("I5").Value = 0.950473488836183
$ws.Range("J5").Value = 0.9504734888361831
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.248835333333334
$ws.Range("N5").Value = 6.746506
$ws.Range("O5").Value = 0.03590294220158827
$ws.Range("P5").Value = 0.03590294220158827
$ws.Range("Q5").Value = 678.1884939828121
$ws.Range("R5").Value = 6103.696445845309
$ws.Range("S5").Value = 0.03412479473382744
$ws.Range("T5").Value = 0.03412479473382744

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mmp2"
$ws.Range("C6").Value = "Sdc2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 301.5732116666666
$ws.Range("H6").Value = 904.7196349999999
$ws.Range("I6").Value = 0.950473488836183
$ws.Range("J6").Value = 0.9504734888361831
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 44.29005966666667
$ws.Range("N6").Value = 132.870179
$ws.Range("O6").Value = 0.7070964373190639
$ws.Range("P6").Value = 0.7070964373190639
$ws.Range("Q6").Value = 13356.69553858496
$ws.Range("R6").Value = 120210.2598472647
$ws.Range("S6").Value = 0.6720764177222861
$ws.Range("T6").Value = 0.6720764177222861

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mmp2"
$ws.Range("C7").Value = "Sdc2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 301.5732116666666
$ws.Range("H7").Value = 904.7196349999999
$ws.Range("I7").Value = 0.950473488836183
$ws.Range("J7").Value = 0.9504734888361831
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 16.09762433333333
$ws.Range("N7").Value = 48.292873
$ws.Range("O7").Value = 0.2570006204793478
$ws.Range("P7").Value = 0.2570006204793479
$ws.Range("Q7").Value = 4854.612270406817
$ws.Range("R7").Value = 43691.51043366135
$ws.Range("S7").Value = 0.2442722763800695
$ws.Range("T7").Value = 0.2442722763800696

# Row 8: M1 -> ECs
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Mmp2"
$ws.Range("C8").Value = "Sdc2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1133353333333333
$ws.Range("H8").Value = 0.340006
$ws.Range("I8").Value = 0.000357200923405664
$ws.Range("J8").Value = 0.000357200923405664
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.248835333333334
$ws.Range("N8").Value = 6.746506
$ws.Range("O8").Value = 0.03590294220158827
$ws.Range("P8").Value = 0.03590294220158827
$ws.Range("Q8").Value = 0.2548725021151111
$ws.Range("R8").Value = 2.293852519036
$ws.Range("S8").Value = 0.00001282456410738752
$ws.Range("T8").Value = 0.00001282456410738752

# Row 9: M1 -> FAPs
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Mmp2"
$ws.Range("C9").Value = "Sdc2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1133353333333333
$ws.Range("H9").Value = 0.340006
$ws.Range("I9").Value = 0.000357200923405664
$ws.Range("J9").Value = 0.000357200923405664
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 44.29005966666667
$ws.Range("N9").Value = 132.870179
$ws.Range("O9").Value = 0.7070964373190639
$ws.Range("P9").Value = 0.7070964373190639
$ws.Range("Q9").Value = 5.01962867567489
$ws.Range("R9").Value = 45.17665808107401
$ws.Range("S9").Value = 0.0002525755003472249
$ws.Range("T9").Value = 0.0002525755003472249

# Row 10: M1 -> sCs
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Mmp2"
$ws.Range("C10").Value = "Sdc2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1133353333333333
$ws.Range("H10").Value = 0.340006
$ws.Range("I10").Value = 0.000357200923405664
$ws.Range("J10").Value = 0.000357200923405664
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 16.09762433333333
$ws.Range("N10").Value = 48.292873
$ws.Range("O10").Value = 0.2570006204793478
$ws.Range("P10").Value = 0.2570006204793479
$ws.Range("Q10").Value = 1.824429619693111
$ws.Range("R10").Value = 16.419866577238
$ws.Range("S10").Value = 0.00009180085895105164
$ws.Range("T10").Value = 0.00009180085895105165

# Row 11: M2 -> ECs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Mmp2"
$ws.Range("C11").Value = "Sdc2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8004729999999999
$ws.Range("H11").Value = 2.401419
$ws.Range("I11").Value = 0.00252286455028413
$ws.Range("J11").Value = 0.002522864550284131
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.248835333333334
$ws.Range("N11").Value = 6.746506
$ws.Range("O11").Value = 0.03590294220158827
$ws.Range("P11").Value = 0.03590294220158827
$ws.Range("Q11").Value = 1.800131965779333
$ws.Range("R11").Value = 16.201187692014
$ws.Range("S11").Value = 0.00009057826013128713
$ws.Range("T11").Value = 0.00009057826013128714

# Row 12: M2 -> FAPs
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Mmp2"
$ws.Range("C12").Value = "Sdc2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.8004729999999999
$ws.Range("H12").Value = 2.401419
$ws.Range("I12").Value = 0.00252286455028413
$ws.Range("J12").Value = 0.002522864550284131
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 44.29005966666667
$ws.Range("N12").Value = 132.870179
$ws.Range("O12").Value = 0.7070964373190639
$ws.Range("P12").Value = 0.7070964373190639
$ws.Range("Q12").Value = 35.45299693155567
$ws.Range("R12").Value = 319.076972384001
$ws.Range("S12").Value = 0.001783908535344471
$ws.Range("T12").Value = 0.001783908535344471

# Row 13: M2 -> sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Mmp2"
$ws.Range("C13").Value = "Sdc2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.8004729999999999
$ws.Range("H13").Value = 2.401419
$ws.Range("I13").Value = 0.00252286455028413
$ws.Range("J13").Value = 0.002522864550284131
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 16.09762433333333
$ws.Range("N13").Value = 48.292873
$ws.Range("O13").Value = 0.2570006204793478
$ws.Range("P13").Value = 0.2570006204793479
$ws.Range("Q13").Value = 12.88571364297633
$ws.Range("R13").Value = 115.971422786787
$ws.Range("S13").Value = 0.0006483777548083723
$ws.Range("T13").Value = 0.0006483777548083725

# Row 14: sCs -> ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Mmp2"
$ws.Range("C14").Value = "Sdc2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.466604
$ws.Range("H14").Value = 7.399812
$ws.Range("I14").Value = 0.007774038338818471
$ws.Range("J14").Value = 0.007774038338818472
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.248835333333334
$ws.Range("N14").Value = 6.746506
$ws.Range("O14").Value = 0.03590294220158827
$ws.Range("P14").Value = 0.03590294220158827
$ws.Range("Q14").Value = 5.546986228541333
$ws.Range("R14").Value = 49.922876056872
$ws.Range("S14").Value = 0.0002791108491515309
$ws.Range("T14").Value = 0.0002791108491515309

# Row 15: sCs -> FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Mmp2"
$ws.Range("C15").Value = "Sdc2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.466604
$ws.Range("H15").Value = 7.399812
$ws.Range("I15").Value = 0.007774038338818471
$ws.Range("J15").Value = 0.007774038338818472
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 44.29005966666667
$ws.Range("N15").Value = 132.870179
$ws.Range("O15").Value = 0.7070964373190639
$ws.Range("P15").Value = 0.7070964373190639
$ws.Range("Q15").Value = 109.2460383340387
$ws.Range("R15").Value = 983.214345006348
$ws.Range("S15").Value = 0.005496994812960355
$ws.Range("T15").Value = 0.005496994812960355

# Row 16: sCs -> sCs
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Mmp2"
$ws.Range("C16").Value = "Sdc2"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.466604
$ws.Range("H16").Value = 7.399812
$ws.Range("I16").Value = 0.007774038338818471
$ws.Range("J16").Value = 0.007774038338818472
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 16.09762433333333
$ws.Range("N16").Value = 48.292873
$ws.Range("O16").Value = 0.2570006204793478
$ws.Range("P16").Value = 0.2570006204793479
$ws.Range("Q16").Value = 39.70646457109733
$ws.Range("R16").Value = 357.358181139876
$ws.Range("S16").Value = 0.001997932676706585
$ws.Range("T16").Value = 0.001997932676706586

